$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set cells that introduce brand-new shared strings first, in the order the
# strings must land in the sharedStrings table:
#   39 "Vukušić, Žagar"  -> first introduced at C7
#   40 "Katarina, Zoki"  -> first introduced at C16
#   41 "nadopuniti"      -> first introduced at F22
#   42 "dodani uvod i opis modela i contollera" -> first introduced at F7
$ws.Range("C7").Value = "Vukušić, Žagar"
$ws.Range("C16").Value = "Katarina, Zoki"
$ws.Range("F22").Value = "nadopuniti"
$ws.Range("F7").Value = "dodani uvod i opis modela i contollera"

# Remaining cells reuse already-existing shared strings, order no longer matters.
$ws.Range("D12").Value = "DA"

$ws.Range("C13").Value = "Vukušić, Žagar"
$ws.Range("D13").Value = "DA"

$ws.Range("C14").Value = "Vukušić, Žagar"
$ws.Range("D14").Value = "DA"

$ws.Range("D15").Value = "NE"

$ws.Range("D16").Value = "DA"

$ws.Range("C17").Value = "Ante"
$ws.Range("D17").Value = "DA"

$ws.Range("D18").Value = "NE"

$ws.Range("C20").Value = "Ante"
$ws.Range("D20").Value = "DA"

$ws.Range("D21").Value = "NE"

$ws.Range("C22").Value = "Ante"
$ws.Range("D22").Value = "DA"

# Match the updated view state: scrolled so row 4 is at top, with C8 selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("C8").Select()
